$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p127v_a1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p127v_1</id>", 2)
